# Automatic update of files.
# - Column C ("Förändrad") on rows 2-42 changes from 45708 -> 45709 (one day later)
# - The last data row (row 43, "A 8130-2025") is removed entirely
# - dimension / used range shrinks accordingly (handled automatically by the
#   engine once the row is deleted)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the "Förändrad" date for every remaining data row (2 through 42).
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 3).Value = 45709
}

# Remove the now-obsolete last row (previously row 43: "A 8130-2025").
$ws.Rows.Item(43).Delete()
